# Sprint_3.xlsx update — "added sprint review protocol"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint review date moved two weeks later (C5) ---
$ws.Range("C5").Value = 45979.8125

# --- Requirement review texts ---------------------------------------------
# NOTE: the order in which brand-new string values are first written controls
# the order they land in the shared-string table, so these are written in a
# specific sequence: B16 text (index 28, reused), then the new "x" marker
# (index 29), then B18's text (index 30) and finally B17's text (index 31).
$ws.Range("B16").Value = "Als Benutzer möchte ich ein Startmenü sehen, damit ich auswählen kann, ob ich lernen, trainieren oder ein Quiz starten will. Als Benutzer möchte ich im Startmenü einfache Navigation haben, damit ich schnell zur gewünschten Funktion gelange. Akzeptanzkriterien: Startmenü existiert und es gibt 3 Auswahlmöglichkeiten: Segmentiertes Herz, 2D Ansicht des Herzens, Quizmodus; nach Auswahl wird entsprechende Szene geladen"
$ws.Range("F16").Value = "x"
$ws.Range("B18").Value = "Als Benutzer möche ich die Möglichkeit haben, zum Startmenü zu navigieren, damit ich den Spielmodus ändern kann. Akzeptanzkriterium: Ein Button für das Startmenü existiert und leitet zum Startmenü weiter. Er kann per Handtracking aktiviert werden."
$ws.Range("B17").Value = "Multiple Choice Quiz implementieren: Als Lernender möchte ich Multiple-Choice-Fragen zu Herzsegmenten beantworten können, damit ich mein Wissen testen kann. Akzeptanzkriterien: Frage und Antwortmöglichkeiten werden in VR dargestellt, Antwort(en) können mittels Handtracking ausgewählt werden, Ausgewählte Antwort wird als richtig/falsch erkannt"

# --- Row 16 (Req 1: Startmenu) ---
$ws.Range("D16").Value = 12
$ws.Rows.Item(16).RowHeight = 126

# --- Row 17 (Req 2: Multiple Choice Quiz) ---
$ws.Range("D17").Value = 23
$ws.Range("F17").Value = "x"
$ws.Rows.Item(17).RowHeight = 94.5

# --- Row 18 (Req 3: Navigate back to start menu) ---
$ws.Range("D18").Value = 8
$ws.Range("F18").Value = "x"

# --- Sheet view: selection moved from B17 to C6:G6, no frozen/scrolled top-left cell ---
$ws.Range("C6:G6").Select() | Out-Null
